$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.118.01"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "1.653.45"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5217"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").Value = "1.649.64"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").Value = "1.878.68"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "0.0₅8011"
$ws.Range("E16").Value = "  -2.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "26.110.48"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.634"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.957"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1205"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.180"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.484"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05704"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.26%  "

$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.361"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.02%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9511"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5667"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01588"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.960"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("D41").Value = "1.059.40"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8423"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("D45").Value = "1.788.91"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05389"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.35%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4394"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.976"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "

